$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared string used by E2 (and soon E3) from "Diego" to "Diego e Iván"
$ws.Range("E2").Value = "Diego e Iván"

# Add the new "paso" entries to row 3
$ws.Range("C3").Value = "ok"
$ws.Range("E3").Value = "Diego e Iván"

# Update the active selection to E2
$ws.Range("E2").Select()
